$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "time_taken" column, matching the style of the
# other header cells in row 1 (bold header style, style index 1).
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Per-row "time_taken" timestamp values (plain, unstyled text cells).
$timestamps = @(
    "2021-10-05 13:40:32.688162",
    "2021-10-05 13:40:32.688173",
    "2021-10-05 13:40:32.688177",
    "2021-10-05 13:40:32.688180",
    "2021-10-05 13:40:32.688182",
    "2021-10-05 13:40:32.688185",
    "2021-10-05 13:40:32.688188",
    "2021-10-05 13:40:32.688190",
    "2021-10-05 13:40:32.688193",
    "2021-10-05 13:40:32.688196",
    "2021-10-05 13:40:32.688198",
    "2021-10-05 13:40:32.688201",
    "2021-10-05 13:40:32.688204",
    "2021-10-05 13:40:32.688207",
    "2021-10-05 13:40:32.688209",
    "2021-10-05 13:40:32.688212",
    "2021-10-05 13:40:32.688215",
    "2021-10-05 13:40:32.688218"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
